$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Marke:text; Farbe:text; Baujahr:Number"
$ws.Range("D3").Value = "Ticketart:select(Bus,Zug,U-Bahn);Häufigkeit:select(Täglich,Wöchentlich,Selten)"
$ws.Range("D4").Value = "Tierart:text; Anzahl:Number"
$ws.Range("D5").Value = "Beschreibung:text"
$ws.Range("D6").Value = "Land:text;Häufigkeit:dropdown(Selten,Oft,Regelmäßig)"
$ws.Range("D7").Value = "Anteil:number;Technik:dropdown(Laptop,PC,Tablet)"
$ws.Range("D8").Value = "Modell:text; Jahr:number"
$ws.Range("D9").Value = "Dringend:checkbox"

$ws.Range("D5").Select()
